$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append row 59 with the new Post49 entry.
# Write F, then E, then C so new shared strings land in the same order
# as the target file (dev.to link, hashnode link, title).
$ws.Range("B59").Value = 49
$ws.Range("F59").Value = "https://dev.to/rahulmishra05/resource-allocation-graph-in-deadlock-operating-system-m04-p02-568p"
$ws.Range("E59").Value = "https://programmingport.hashnode.dev/resource-allocation-graph-in-deadlock-or-operating-system-m04-p02"
$ws.Range("C59").Value = "Resource Allocation Graph in Deadlock | Operating System - M04 P02"
$ws.Range("D59").Value = [DateTime]"2020-12-08"

# Match the look of the rest of the table: hyperlink-style link columns.
$ws.Range("E59").Style = "Hyperlink"
$ws.Range("F59").Style = "Hyperlink"

# Extend the table range to include the new row.
$table = $ws.ListObjects.Item("Table2")
$table.Resize($ws.Range("B10:F59"))

$ws.Range("E59").Select()
